# Daily attendance processing - 2025-11-04 17:45:48
# Normalizes the "Recorded By" list in column G for each attendance row by
# rotating the comma-separated list of recorders left by one position
# (moving the first contributor to the end of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value()

    if ($null -eq $raw) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    if ($parts.Count -le 1) {
        continue
    }

    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $rotated = $trimmed[1..($trimmed.Count - 1)] + $trimmed[0]
    $newText = $rotated -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
